$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2741.1
$ws.Range("I19").Value = 4043.5
$ws.Range("J19").Value = 787.5
$ws.Range("K19").Value = 4043.5
$ws.Range("L19").Value = 787.5
$ws.Range("M19").Value = -3868.5
$ws.Range("N19").Value = -1137.5
$ws.Range("H112").Value = 3704782.2
$ws.Range("J112").Value = 3704782.2
$ws.Range("L112").Value = 11114346.6
$ws.Range("N112").Value = -11116562.6
$ws.Range("H137").Value = 2135.2083
$ws.Range("I137").Value = 2046.9445
$ws.Range("J137").Value = 2400
$ws.Range("K137").Value = 6140.833500000001
$ws.Range("L137").Value = 7200
$ws.Range("M137").Value = -3590.833500000001
$ws.Range("N137").Value = -12300
$ws.Range("H141").Value = 1426.6
$ws.Range("I141").Value = 1007.4
$ws.Range("J141").Value = 3941.8
$ws.Range("K141").Value = 3022.2
$ws.Range("L141").Value = 11825.4
$ws.Range("M141").Value = 2157.8
$ws.Range("N141").Value = -22185.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1765.56
$ws.Range("I32").Value = 1777.5604
$ws.Range("K32").Value = 1777.5604
$ws.Range("M32").Value = -1490.5604
$ws.Range("H61").Value = 515466.22
$ws.Range("I61").Value = 546645.5600000001
$ws.Range("J61").Value = 1007
$ws.Range("K61").Value = 546645.5600000001
$ws.Range("L61").Value = 1007
$ws.Range("M61").Value = -546433.5600000001
$ws.Range("N61").Value = -1431
$ws.Range("H74").Value = 2263.8518
$ws.Range("I74").Value = 2263.8518
$ws.Range("K74").Value = 2263.8518
$ws.Range("M74").Value = -1389.8518
$ws.Range("H77").Value = 2263.8518
$ws.Range("I77").Value = 2263.8518
$ws.Range("K77").Value = 11319.259
$ws.Range("M77").Value = -6951.259
$ws.Range("H132").Value = 12575
$ws.Range("I132").Value = 1409.2565
$ws.Range("K132").Value = 4227.7695
$ws.Range("M132").Value = -1697.7695
$ws.Range("H136").Value = 515466.22
$ws.Range("I136").Value = 546645.5600000001
$ws.Range("J136").Value = 1007
$ws.Range("K136").Value = 1639936.68
$ws.Range("L136").Value = 3021
$ws.Range("M136").Value = -1637386.68
$ws.Range("N136").Value = -8121
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3102.6304
$ws.Range("I134").Value = 3381.139
$ws.Range("J134").Value = 2100
$ws.Range("K134").Value = 10143.417
$ws.Range("L134").Value = 6300
$ws.Range("M134").Value = -7608.417000000001
$ws.Range("N134").Value = -11370
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4375.8887
$ws.Range("I31").Value = 2811.4285
$ws.Range("J31").Value = 6060.6924
$ws.Range("K31").Value = 2811.4285
$ws.Range("L31").Value = 6060.6924
$ws.Range("M31").Value = -2516.4285
$ws.Range("N31").Value = -6650.6924
$ws.Range("H34").Value = 4375.8887
$ws.Range("I34").Value = 2811.4285
$ws.Range("J34").Value = 6060.6924
$ws.Range("K34").Value = 2811.4285
$ws.Range("L34").Value = 6060.6924
$ws.Range("M34").Value = -2609.4285
$ws.Range("N34").Value = -6464.6924
$ws.Range("H99").Value = 21742866
$ws.Range("I99").Value = 3550
$ws.Range("J99").Value = 38465416
$ws.Range("K99").Value = 3550
$ws.Range("L99").Value = 38465416
$ws.Range("M99").Value = -2052
$ws.Range("N99").Value = -38468412
$ws.Range("H122").Value = 2857.2856
$ws.Range("I122").Value = 2816.8333
$ws.Range("J122").Value = 3100
$ws.Range("K122").Value = 8450.499899999999
$ws.Range("L122").Value = 9300
$ws.Range("M122").Value = -6000.499899999999
$ws.Range("N122").Value = -14200
$ws.Range("H126").Value = 21742866
$ws.Range("I126").Value = 3550
$ws.Range("J126").Value = 38465416
$ws.Range("K126").Value = 10650
$ws.Range("L126").Value = 115396248
$ws.Range("M126").Value = -8180
$ws.Range("N126").Value = -115401188
$ws.Range("H134").Value = 808.54
$ws.Range("I134").Value = 732.5106
$ws.Range("J134").Value = 1999.6666
$ws.Range("K134").Value = 2197.5318
$ws.Range("L134").Value = 5998.9998
$ws.Range("M134").Value = 337.4682000000003
$ws.Range("N134").Value = -11068.9998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 6764.3
$ws.Range("I62").Value = 3592
$ws.Range("K62").Value = 10776
$ws.Range("M62").Value = -10090
$ws.Range("H65").Value = 6764.3
$ws.Range("I65").Value = 3592
$ws.Range("K65").Value = 32328
$ws.Range("M65").Value = -28896
$ws.Range("H75").Value = 814.1429000000001
$ws.Range("I75").Value = 400
$ws.Range("J75").Value = 979.8
$ws.Range("K75").Value = 1200
$ws.Range("L75").Value = 2939.4
$ws.Range("M75").Value = -202
$ws.Range("N75").Value = -4935.4
$ws.Range("H78").Value = 814.1429000000001
$ws.Range("I78").Value = 400
$ws.Range("J78").Value = 979.8
$ws.Range("K78").Value = 3600
$ws.Range("L78").Value = 8818.199999999999
$ws.Range("M78").Value = 1392
$ws.Range("N78").Value = -18802.2
$ws.Range("H113").Value = 671.46155
$ws.Range("I113").Value = 506.25
$ws.Range("K113").Value = 1518.75
$ws.Range("M113").Value = 651.25
$ws.Range("H127").Value = 964
$ws.Range("J127").Value = 964
$ws.Range("L127").Value = 2892
$ws.Range("N127").Value = -12812
$ws.Range("H131").Value = 718.12
$ws.Range("J131").Value = 735.23914
$ws.Range("L131").Value = 2205.71742
$ws.Range("N131").Value = -12285.71742
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 78432824
$ws.Range("I122").Value = 25641928
$ws.Range("J122").Value = 250003230
$ws.Range("K122").Value = 76925784
$ws.Range("L122").Value = 750009690
$ws.Range("M122").Value = -76923334
$ws.Range("N122").Value = -750014590
$ws.Range("H132").Value = 31614.5
$ws.Range("I132").Value = 4437.067
$ws.Range("J132").Value = 167501.67
$ws.Range("K132").Value = 13311.201
$ws.Range("L132").Value = 502505.01
$ws.Range("M132").Value = -10781.201
$ws.Range("N132").Value = -507565.01
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5112.75
$ws.Range("I22").Value = 3483.6667
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 3483.6667
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -3188.6667
$ws.Range("N22").Value = -10590
$ws.Range("H27").Value = 5112.75
$ws.Range("I27").Value = 3483.6667
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 3483.6667
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = -3376.6667
$ws.Range("N27").Value = -10214
$ws.Range("H132").Value = 1187.0333
$ws.Range("I132").Value = 1210.7241
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 3632.1723
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -1102.1723
$ws.Range("N132").Value = -6560
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1351936
$ws.Range("I113").Value = 745.6
$ws.Range("K113").Value = 2236.8
$ws.Range("M113").Value = -66.80000000000018
